$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.113.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.607.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.86"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.16"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.598.24"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.55%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.67"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +20.93%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.16"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000283"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.186.49"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "667.12"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.84"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.608.66"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.129.94"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.60%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.66"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.31"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.926"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.40"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.66%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.74%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.91"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.39"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.33"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.93"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.37"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.92"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "576.87"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.34%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.89"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.22%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.566.03"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.69%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.81%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.342"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "34.46"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.88%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.76%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.15%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.06%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.86"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.69%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.89"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.54%  "
